$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteSpecial paste-type constants
$xlPasteValues  = -4163
$xlPasteFormats = -4122

# --- Swap the D:G contents (values + formats) between row 2 and row 3 ---
# Stage row 2's original D:G block in a scratch range first (values, then
# formats, so each PasteSpecial gets a fresh, un-mutated clipboard source).

$ws.Range("D2:G2").Copy()
$ws.Range("Z1:AC1").PasteSpecial($xlPasteValues)
$ws.Range("D2:G2").Copy()
$ws.Range("Z1:AC1").PasteSpecial($xlPasteFormats)

# Move row 3's D:G block onto row 2.
$ws.Range("D3:G3").Copy()
$ws.Range("D2:G2").PasteSpecial($xlPasteValues)
$ws.Range("D3:G3").Copy()
$ws.Range("D2:G2").PasteSpecial($xlPasteFormats)

# Move the staged (original row 2) D:G block onto row 3.
$ws.Range("Z1:AC1").Copy()
$ws.Range("D3:G3").PasteSpecial($xlPasteValues)
$ws.Range("Z1:AC1").Copy()
$ws.Range("D3:G3").PasteSpecial($xlPasteFormats)

# Drop the scratch range so it doesn't linger in the saved sheet.
$ws.Range("Z1:AC1").Clear()

# D2 ends up re-keyed with a quote-prefixed, centered numeric style (the
# same one already used by G2/G3) rather than the plain style D3 used to
# carry - line up D2's format with G2's to match.
$ws.Range("G2").Copy()
$ws.Range("D2").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# --- Move the sheet's active selection ---
$ws.Range("J5").Select()
